$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$aw = $excel.ActiveWindow
$aw.ScrollRow = 6
$aw.ScrollColumn = 3
